$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# like "1.003" or "85.00" keep their exact text representation
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '30.255.79'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '1.855.22'
$ws.Range('E3').Value = '  -1.09%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').Value = '233.19'
$ws.Range('E5').Value = '  -2.32%  '
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('D7').Value = '0.4750'
$ws.Range('E7').Value = '  -0.97%  '
$ws.Range('D8').Value = '0.2769'
$ws.Range('E8').Value = '  -2.18%  '
$ws.Range('D9').Value = '0.06429'
$ws.Range('E9').Value = '  -1.64%  '
$ws.Range('D10').Value = '1.856.96'
$ws.Range('E10').Value = '  -1.01%  '
$ws.Range('D11').Value = '0.07413'
$ws.Range('E11').Value = '  -0.72%  '
$ws.Range('D12').Value = '15.99'
$ws.Range('E12').Value = '  -4.37%  '
$ws.Range('D13').Value = '4.995'
$ws.Range('E13').Value = '  -2.13%  '
$ws.Range('D14').Value = '85.00'
$ws.Range('E14').Value = '  -3.80%  '
$ws.Range('D15').Value = '0.6321'
$ws.Range('E15').Value = '  -4.73%  '
$ws.Range('D16').Value = '30.241.67'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').Value = '1.002'
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('D18').Value = '12.78'
$ws.Range('E18').Value = '  -4.36%  '
$ws.Range('D19').Value = '0.000007324'
$ws.Range('E19').Value = '  -4.03%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '223.00'
$ws.Range('E20').Value = '  +1.80%  '
$ws.Range('B21').Value = 'BinanceUSD'
$ws.Range('C21').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D21').Value = '1.003'
$ws.Range('E21').Value = '  +0.38%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '5.084'
$ws.Range('E22').Value = '  -4.37%  '
$ws.Range('B23').Value = 'BitDAO'
$ws.Range('C23').Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range('D23').Value = '0.3859'
$ws.Range('E23').Value = '  -4.18%  '
$ws.Range('D24').Value = '5.998'
$ws.Range('E24').Value = '  -3.90%  '
$ws.Range('D25').Value = '9.215'
$ws.Range('E25').Value = '  -1.54%  '
$ws.Range('D26').Value = '166.56'
$ws.Range('E26').Value = '  -0.44%  '
$ws.Range('D27').Value = '17.71'
$ws.Range('E27').Value = '  -4.08%  '
$ws.Range('D28').Value = '1.865'
$ws.Range('E28').Value = '  -5.94%  '
$ws.Range('D29').Value = '0.1026'
$ws.Range('E29').Value = '  +8.75%  '
$ws.Range('D30').Value = '1.375'
$ws.Range('E30').Value = '  -5.72%  '
$ws.Range('D31').Value = '4.209'
$ws.Range('E31').Value = '  -2.47%  '
$ws.Range('D32').Value = '3.890'
$ws.Range('E32').Value = '  -3.93%  '
$ws.Range('D33').Value = '0.04888'
$ws.Range('E33').Value = '  -3.72%  '
$ws.Range('D34').Value = '1.150'
$ws.Range('E34').Value = '  -4.80%  '
$ws.Range('D35').Value = '0.7270'
$ws.Range('E35').Value = '  -3.51%  '
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').Value = '  +0.33%  '
$ws.Range('D37').Value = '2.713'
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('D38').Value = '0.01901'
$ws.Range('E38').Value = '  +3.91%  '
$ws.Range('D39').Value = '2.621'
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('D40').Value = '0.9021'
$ws.Range('E40').Value = '  -0.75%  '
$ws.Range('D41').Value = '1.981'
$ws.Range('E41').Value = '  -4.66%  '
$ws.Range('D42').Value = '105.32'
$ws.Range('E42').Value = '  -1.57%  '
$ws.Range('D43').Value = '0.9945'
$ws.Range('E43').Value = '  -0.95%  '
$ws.Range('D44').Value = '0.4096'
$ws.Range('E44').Value = '  -4.67%  '
$ws.Range('D45').Value = '5.550'
$ws.Range('E45').Value = '  -6.09%  '
$ws.Range('D46').Value = '7.042'
$ws.Range('E46').Value = '  -5.62%  '
$ws.Range('D47').Value = '61.10'
$ws.Range('E47').Value = '  -5.61%  '
$ws.Range('D48').Value = '0.1205'
$ws.Range('E48').Value = '  -6.23%  '
$ws.Range('D49').Value = '8.774'
$ws.Range('E49').Value = '  -2.05%  '
$ws.Range('D50').Value = '1.397'
$ws.Range('E50').Value = '  -5.68%  '
$ws.Range('D51').Value = '0.05597'
$ws.Range('E51').Value = '  -0.67%  '
